# TRE-210-BE: add an accumulated-points column ("Tổng điểm tích lũy") to the
# revenue export template, right before the existing "Tổng doanh thu" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at H; Excel shifts the old H:L columns to I:M and
# copies the formatting of the column immediately to the left (G) onto the
# newly inserted column - including the row-8 header cell style.
$ws.Columns("H:H").Insert()

# Give the new column (H) the same width as its neighbour (G) so that it
# keeps blending into the F:H run of equally-wide columns.
$ws.Columns("H:H").ColumnWidth = $ws.Columns("G:G").ColumnWidth

# Header text for the newly inserted column.
$ws.Cells.Item(8, 8).Value = "Tổng điểm tích lũy"

# Move the active selection to G14, matching the saved view state.
$ws.Range("G14").Select()
